$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '60.290.40'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -0.07%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.334.18'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -0.55%  '

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.01%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '548.25'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +0.20%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '131.37'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -0.85%  '

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.03%  '

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -1.35%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.333.33'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -0.46%  '

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +0.70%  '

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +1.29%  '

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -0.67%  '

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +0.68%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '23.70'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -0.91%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.750.40'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -0.44%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '60.249.65'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -0.08%  '

$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +0.79%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.333.61'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -1.37%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.66'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +0.12%  '

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -1.63%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '314.03'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -0.13%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.58'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -3.66%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.00'

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '64.09'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +0.89%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.170'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -1.32%  '

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +0.05%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.92'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +0.18%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.38'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +1.43%  '

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +7.58%  '

$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -0.74%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '171.10'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -0.27%  '

$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +0.12%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.10'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +2.46%  '

$ws.Range("B34").Value = 'ImmutableX'
$ws.Range("C34").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.37'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -2.82%  '

$ws.Range("B35").Value = 'PolygonEcosystemToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.384'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +0.56%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '18.08'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -0.05%  '

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +0.00%  '

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -1.95%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '323.87'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -0.61%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '38.13'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +0.09%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.53'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -0.10%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '137.55'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -2.78%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.50'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +1.10%  '

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +0.34%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '19.30'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -1.51%  '

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +0.91%  '

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -0.02%  '

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +1.20%  '

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +4.61%  '

$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -0.90%  '
